# Add columns I (I0) and J (IF) to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - mirror style of existing header cells (bold font,
# thin box border, centered horizontally, top-aligned vertically)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous
$headerRange.Borders.Weight = 2           # xlThin

# Data values for I2:J27
$values = @(
    @(5, 5),
    @(6, 6),
    @(8, 8),
    @(6, 9),
    @(5, 6),
    @(5, 6),
    @(1, 4),
    @(1, 4),
    @(5, 8),
    @(5, 8),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(1, 6),
    @(5, 9),
    @(1, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
